# Add links to datasheets
#
# The BOM sheet gets a new "dsheet" column inserted between the existing
# "URL" (C) and "Cost per unit" (D, which shifts to E) columns. Two of the
# rows get a link to the part's datasheet PDF; the rest are left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Insert a new blank column at D - this shifts the old D:G (Cost per unit,
# QTY per device, Postage cost, Cost) out to E:H, and (like a native Excel
# "Insert" of an entire column) carries the formatting of column C into the
# new column D for each row.
$ws.Columns("D:D").Insert()

# Header for the new column.
$ws.Range("D1").Value = "dsheet"

# Datasheet links for the rows that have one.
$ws.Range("D5").Value = "https://au.mouser.com/datasheet/2/389/stm32f101rc-956301.pdf"
$ws.Range("D3").Value = "https://cdn-shop.adafruit.com/datasheets/62684.pdf"

# Leave the selection on the newly added cell, matching the author's
# last-touched location.
$ws.Range("D3").Select()
